$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.07686466666666666
$ws.Range("H2").Value = 0.230594
$ws.Range("I2").Value = 0.07882978817959985
$ws.Range("J2").Value = 0.07882978817959985
$ws.Range("M2").Value = 0.24449
$ws.Range("N2").Value = 0.73347
$ws.Range("O2").Value = 0.009675524511058336
$ws.Range("P2").Value = 0.009675524511058334
$ws.Range("Q2").Value = 0.01879264235333333
$ws.Range("R2").Value = 0.16913378118
$ws.Range("S2").Value = 0.0007627195477332549
$ws.Range("T2").Value = 0.0007627195477332548
$ws.Range("G3").Value = 0.07686466666666666
$ws.Range("H3").Value = 0.230594
$ws.Range("I3").Value = 0.07882978817959985
$ws.Range("J3").Value = 0.07882978817959985
$ws.Range("O3").Value = 0.1748614367985708
$ws.Range("P3").Value = 0.1748614367985708
$ws.Range("Q3").Value = 0.3396310390604445
$ws.Range("R3").Value = 3.056679351544
$ws.Range("S3").Value = 0.01378429002361182
$ws.Range("T3").Value = 0.01378429002361182
$ws.Range("G4").Value = 0.07686466666666666
$ws.Range("H4").Value = 0.230594
$ws.Range("I4").Value = 0.07882978817959985
$ws.Range("J4").Value = 0.07882978817959985
$ws.Range("M4").Value = 20.60586566666667
$ws.Range("N4").Value = 61.81759700000001
$ws.Range("O4").Value = 0.815463038690371
$ws.Range("P4").Value = 0.8154630386903708
$ws.Range("Q4").Value = 1.583862995846445
$ws.Range("R4").Value = 14.254766962618
$ws.Range("S4").Value = 0.06428277860825478
$ws.Range("T4").Value = 0.06428277860825476
$ws.Range("G5").Value = 0.3352123333333334
$ws.Range("I5").Value = 0.3437823694266471
$ws.Range("J5").Value = 0.3437823694266471
$ws.Range("M5").Value = 0.24449
$ws.Range("N5").Value = 0.73347
$ws.Range("O5").Value = 0.009675524511058336
$ws.Range("P5").Value = 0.009675524511058334
$ws.Range("Q5").Value = 0.08195606337666668
$ws.Range("R5").Value = 0.73760457039
$ws.Range("S5").Value = 0.003326274741857236
$ws.Range("T5").Value = 0.003326274741857236
$ws.Range("G6").Value = 0.3352123333333334
$ws.Range("I6").Value = 0.3437823694266471
$ws.Range("J6").Value = 0.3437823694266471
$ws.Range("O6").Value = 0.1748614367985708
$ws.Range("P6").Value = 0.1748614367985708
$ws.Range("Q6").Value = 1.481155360623556
$ws.Range("S6").Value = 0.06011427906396058
$ws.Range("T6").Value = 0.06011427906396057
$ws.Range("G7").Value = 0.3352123333333334
$ws.Range("I7").Value = 0.3437823694266471
$ws.Range("J7").Value = 0.3437823694266471
$ws.Range("M7").Value = 20.60586566666667
$ws.Range("N7").Value = 61.81759700000001
$ws.Range("O7").Value = 0.815463038690371
$ws.Range("P7").Value = 0.8154630386903708
$ws.Range("Q7").Value = 6.907340310476558
$ws.Range("R7").Value = 62.16606279428901
$ws.Range("S7").Value = 0.2803418156208294
$ws.Range("T7").Value = 0.2803418156208293
$ws.Range("G8").Value = 0.5629943333333333
$ws.Range("H8").Value = 1.688983
$ws.Range("I8").Value = 0.5773878423937531
$ws.Range("J8").Value = 0.5773878423937531
$ws.Range("M8").Value = 0.24449
$ws.Range("N8").Value = 0.73347
$ws.Range("O8").Value = 0.009675524511058336
$ws.Range("P8").Value = 0.009675524511058334
$ws.Range("Q8").Value = 0.1376464845566666
$ws.Range("R8").Value = 1.23881836101
$ws.Range("S8").Value = 0.005586530221467845
$ws.Range("T8").Value = 0.005586530221467844
$ws.Range("G9").Value = 0.5629943333333333
$ws.Range("H9").Value = 1.688983
$ws.Range("I9").Value = 0.5773878423937531
$ws.Range("J9").Value = 0.5773878423937531
$ws.Range("O9").Value = 0.1748614367985708
$ws.Range("P9").Value = 0.1748614367985708
$ws.Range("Q9").Value = 2.487623490834222
$ws.Range("R9").Value = 22.388611417508
$ws.Range("S9").Value = 0.1009628677109984
$ws.Range("T9").Value = 0.1009628677109984
$ws.Range("G10").Value = 0.5629943333333333
$ws.Range("H10").Value = 1.688983
$ws.Range("I10").Value = 0.5773878423937531
$ws.Range("J10").Value = 0.5773878423937531
$ws.Range("M10").Value = 20.60586566666667
$ws.Range("N10").Value = 61.81759700000001
$ws.Range("O10").Value = 0.815463038690371
$ws.Range("P10").Value = 0.8154630386903708
$ws.Range("Q10").Value = 11.60098560376122
$ws.Range("R10").Value = 104.408870433851
$ws.Range("S10").Value = 0.4708384444612869
$ws.Range("T10").Value = 0.4708384444612868
